$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove three TECHNIQUE contacts (column C) that left the company
$ws.Range("C5").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("C7").ClearContents()

# Add a new TECHNIQUE contact in row 19
$ws.Range("C19").Value = "Chérif Hemmet"

# Add a new RH contact for Thomson Video (row 8)
$ws.Range("E8").Value = "Laurence Fabbroni"

# Add a brand new row 22 with a new CLIENT contact
$ws.Range("B22").Value = "Charles-Antoine Robin"

$ws.Range("B17").Select()
